$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Trim the trailing newline from the METAL1.SP.1.1 rule text in C2
$c2 = "rule METAL1.SP.1.1 {`n    caption METAL1.SP.1.1: Metal1 to Metal1 spacing must be >= 0.06 um;`n    exte Metal1 Metal1 -lt 0.06 -output region -singular -abut lt 90;`n}"
$ws.Range("C2").Value = $c2

# Trim the trailing blank lines from the "drc = exte" text in D3
$ws.Range("D3").Value = "drc = exte"

# Introduce an empty D2 cell (style index 0, same as an untouched/default cell)
# so the row spans through column D, matching the new dimension/content.
$ws.Range("D2").Value = 0
$ws.Range("D2").ClearContents()
$ws.Range("D2").WrapText = $false
